$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 27, pushing existing data
# (old rows 27..127) down to rows 29..129. The new rows inherit the
# formatting (incl. date style on column D) from the row above, matching
# how Excel's Rows.Insert works.
$ws.Rows("27:28").Insert()

# Populate the first new row (27) with its data.
$ws.Range("A27").Value = 1
$ws.Range("B27").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C27").Value = "Arica y Parinacota"
$ws.Range("D27").Value = 44998
$ws.Range("E27").Value = 15
$ws.Range("F27").Value = 100112021
$ws.Range("G27").Value = "Ají"
$ws.Range("H27").Value = "Inferno"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 120
$ws.Range("K27").Value = 31000
$ws.Range("L27").Value = 32000
$ws.Range("M27").Value = 31500
$ws.Range("N27").Value = "$/caja 15 kilos"
$ws.Range("O27").Value = "Región de Arica y Parinacota"
$ws.Range("P27").Value = 2100
$ws.Range("Q27").Value = 15
$ws.Range("R27").Value = "Hortaliza"

# Populate the second new row (28) with its data.
$ws.Range("A28").Value = 1
$ws.Range("B28").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C28").Value = "Arica y Parinacota"
$ws.Range("D28").Value = 44998
$ws.Range("E28").Value = 15
$ws.Range("F28").Value = 100112021
$ws.Range("G28").Value = "Ají"
$ws.Range("H28").Value = "Inferno"
$ws.Range("I28").Value = "Segunda"
$ws.Range("J28").Value = 100
$ws.Range("K28").Value = 27000
$ws.Range("L28").Value = 28000
$ws.Range("M28").Value = 27500
$ws.Range("N28").Value = "$/caja 15 kilos"
$ws.Range("O28").Value = "Región de Arica y Parinacota"
$ws.Range("P28").Value = 1833
$ws.Range("Q28").Value = 15
$ws.Range("R28").Value = "Hortaliza"
